$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the commit diff.
# Column D cells whose new value would be auto-parsed as a number by Excel
# are written with a Text number format first so they stay text (matching the
# original inline-string storage of this sheet).

$ws.Range('D2').Value = '69.592.46'
$ws.Range('E2').Value = '  +3.30%  '
$ws.Range('D3').Value = '3.372.83'
$ws.Range('E3').Value = '  +4.72%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '191.62'
$ws.Range('E5').Value = '  +4.37%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '593.09'
$ws.Range('E6').Value = '  +2.55%  '
$ws.Range('B7').Value = 'USDC'
$ws.Range('C7').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('B8').Value = 'XRP'
$ws.Range('C8').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.608'
$ws.Range('E8').Value = '  +0.36%  '
$ws.Range('E9').Value = '  +3.49%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.77'
$ws.Range('E10').Value = '  +3.12%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.421'
$ws.Range('E11').Value = '  +2.78%  '
$ws.Range('D12').Value = '3.960.91'
$ws.Range('E12').Value = '  +4.82%  '
$ws.Range('E13').Value = '  +1.44%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.66'
$ws.Range('D15').Value = '69.616.74'
$ws.Range('E15').Value = '  +3.19%  '
$ws.Range('E16').Value = '  +2.29%  '
$ws.Range('D17').Value = '3.346.94'
$ws.Range('E17').Value = '  +3.76%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '451.81'
$ws.Range('E18').Value = '  +14.21%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.84'
$ws.Range('E19').Value = '  +1.74%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.85'
$ws.Range('E20').Value = '  +3.20%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.81'
$ws.Range('E21').Value = '  +3.83%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '74.72'
$ws.Range('E22').Value = '  +4.95%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.16%  '
$ws.Range('B24').Value = 'WrappedeETH'
$ws.Range('C24').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D24').Value = '3.517.90'
$ws.Range('E24').Value = '  +4.60%  '
$ws.Range('B25').Value = 'PEPE'
$ws.Range('C25').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000123'
$ws.Range('E25').Value = '  +5.36%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.519'
$ws.Range('E26').Value = '  +0.83%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.191'
$ws.Range('E27').Value = '  +4.07%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.54'
$ws.Range('E28').Value = '  +0.52%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.22%  '
$ws.Range('E30').Value = '  +2.47%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '23.32'
$ws.Range('E31').Value = '  +3.41%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.65'
$ws.Range('E32').Value = '  +2.18%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.29'
$ws.Range('E33').Value = '  +3.34%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.02'
$ws.Range('E34').Value = '  +1.09%  '
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.54'
$ws.Range('E36').Value = '  +5.04%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '165.01'
$ws.Range('E37').Value = '  +2.93%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.94'
$ws.Range('E38').Value = '  +3.57%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '27.24'
$ws.Range('E39').Value = '  +3.26%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.819'
$ws.Range('E40').Value = '  +2.21%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.62'
$ws.Range('E41').Value = '  +1.81%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.56'
$ws.Range('E42').Value = '  +0.93%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '2.738.11'
$ws.Range('E43').Value = '  +5.70%  '
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.55'
$ws.Range('E44').Value = '  +3.76%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '25.73'
$ws.Range('E45').Value = '  +5.19%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0691'
$ws.Range('E46').Value = '  +1.34%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '343.24'
$ws.Range('E47').Value = '  +3.17%  '
$ws.Range('E48').Value = '  +0.10%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0286'
$ws.Range('E49').Value = '  +3.29%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '33.02'
$ws.Range('E50').Value = '  +8.43%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.02'
$ws.Range('E51').Value = '  +6.07%  '
